# "Natmi following Dr Hou advice"
# The NATMI Il1a-Il1rap ligand-receptor pair sheet was recomputed with a
# third "Neutro" sending-cluster group added alongside the existing M1/M2
# clusters. This refreshes the existing M1/M2 rows with their recalculated
# statistics and appends the six new Neutro rows (each against the same
# six target clusters: ECs, FAPs, M1, M2, Neutro, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
# D Target cluster | E Ligand-expressing cells | F Ligand detection rate |
# G Ligand average expression value | H Ligand total expression value |
# I Ligand derived specificity (avg) | J Ligand derived specificity (total) |
# K Receptor-expressing cells | L Receptor detection rate |
# M Receptor average expression value | N Receptor total expression value |
# O Receptor derived specificity (avg) | P Receptor derived specificity (total) |
# Q Edge average expression weight | R Edge total expression weight |
# S Edge average expression derived specificity | T Edge total expression derived specificity
$rows = @(
    @("M1", "Il1a", "Il1rap", "ECs", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 2, 1, 1.974591, 3.949182, 0.06180909787364696, 0.04563698887180197, 6.034494241143001, 36.20696544685801, 0.02818194540679362, 0.02080824948367239),
    @("M1", "Il1a", "Il1rap", "FAPs", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 3, 1, 7.979832333333334, 23.939497, 0.2497865318469965, 0.2766462923677706, 24.386950138427, 219.482551245843, 0.1138905217198734, 0.1261372674365544),
    @("M1", "Il1a", "Il1rap", "M1", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 3, 1, 7.758748, 23.276244, 0.2428660996170663, 0.2689817001103893, 23.711300276604, 213.401702489436, 0.1107351408778168, 0.122642585779747),
    @("M1", "Il1a", "Il1rap", "M2", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 3, 1, 4.653158666666667, 13.959476, 0.1456542339398937, 0.1613165589401012, 14.220392565916, 127.983533093244, 0.06641125352700815, 0.07355251271512359),
    @("M1", "Il1a", "Il1rap", "Neutro", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 3, 1, 2.249720666666667, 6.749162, 0.07042126945497384, 0.07799372910339122, 6.875310586942001, 61.877795282478, 0.0321086771936747, 0.03556135085739815),
    @("M1", "Il1a", "Il1rap", "sCs", 3, 1, 3.056073, 9.168219, 0.4559514113020136, 0.4559514113020135, 2, 1, 7.330557000000001, 14.661114, 0.2294627672674229, 0.1694247306065459, 22.402717322661, 134.416303935966, 0.1046238725768469, 0.07724944502951804),
    @("M2", "Il1a", "Il1rap", "ECs", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 2, 1, 1.974591, 3.949182, 0.06180909787364696, 0.04563698887180197, 7.093005340149, 42.55803204089401, 0.03312534261832561, 0.02445822612615257),
    @("M2", "Il1a", "Il1rap", "FAPs", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 3, 1, 7.979832333333334, 23.939497, 0.2497865318469965, 0.2766462923677706, 28.66466693802767, 257.982002442249, 0.1338680668950984, 0.1482630152199497),
    @("M2", "Il1a", "Il1rap", "M1", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 3, 1, 7.758748, 23.276244, 0.2428660996170663, 0.2689817001103893, 27.870501281972, 250.834511537748, 0.1301592004568281, 0.1441553311849143),
    @("M2", "Il1a", "Il1rap", "M2", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 3, 1, 4.653158666666667, 13.959476, 0.1456542339398937, 0.1613165589401012, 16.71479271972133, 150.433134477492, 0.07806045661646614, 0.08645436462806723),
    @("M2", "Il1a", "Il1rap", "Neutro", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 3, 1, 2.249720666666667, 6.749162, 0.07042126945497384, 0.07799372910339122, 8.081309345839333, 72.73178411255401, 0.03774086272998369, 0.04179917014663699),
    @("M2", "Il1a", "Il1rap", "sCs", 3, 1, 3.592139, 10.776417, 0.5359298834298145, 0.5359298834298145, 2, 1, 7.330557000000001, 14.661114, 0.2294627672674229, 0.1694247306065459, 26.332379691423, 157.994278148538, 0.1229759541131126, 0.09079977612409387),
    @("Neutro", "Il1a", "Il1rap", "ECs", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 2, 1, 1.974591, 3.949182, 0.06180909787364696, 0.04563698887180197, 0.10745066025, 0.6447039615000001, 0.0005018098485277302, 0.000370513261977001),
    @("Neutro", "Il1a", "Il1rap", "FAPs", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 3, 1, 7.979832333333334, 23.939497, 0.2497865318469965, 0.2766462923677706, 0.434235876138889, 3.908122885250001, 0.002027943232024598, 0.002246009711266442),
    @("Neutro", "Il1a", "Il1rap", "M1", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 3, 1, 7.758748, 23.276244, 0.2428660996170663, 0.2689817001103893, 0.4222052036666667, 3.799846833, 0.001971758282421437, 0.002183783145728051),
    @("Neutro", "Il1a", "Il1rap", "M2", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 3, 1, 4.653158666666667, 13.959476, 0.1456542339398937, 0.1613165589401012, 0.2532093841111112, 2.278884457, 0.001182523796419357, 0.001309681596910362),
    @("Neutro", "Il1a", "Il1rap", "Neutro", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 3, 1, 2.249720666666667, 6.749162, 0.07042126945497384, 0.07799372910339122, 0.1224222996111111, 1.1018006965, 0.0005717295313154491, 0.0006332080993560743),
    @("Neutro", "Il1a", "Il1rap", "sCs", 1, 0.3333333333333333, 0.05441666666666667, 0.16325, 0.0081187052681719, 0.008118705268171898, 2, 1, 7.330557000000001, 14.661114, 0.2294627672674229, 0.1694247306065459, 0.3989044767500001, 2.3934268605, 0.001862940577463329, 0.001375509452933969)
)

$startRow = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($startRow, $c + 1).Value = $row[$c]
    }
    $startRow = $startRow + 1
}

Write-Host "done"